# Add "2022-Q4" worksheet data, pushing the existing quarterly sheets down
# (总计 stays first; 2022-Q4 is inserted right after it; 2020-Q4 ends up last).

$wb = $excel.ActiveWorkbook

# --- 1) Create a brand-new worksheet (appended at the end for now) -------
$totalSheet = $wb.Worksheets.Item(1)
$lastSheet  = $wb.Worksheets.Item($wb.Worksheets.Count)

$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2022-Q4"

# --- 2) Populate the new "2022-Q4" sheet with the fund table -------------
# (Do this BEFORE moving the sheet - the sheet reference/position tracking
# is only reliable while the sheet stays where it was created.)
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $newSheet.Cells.Item(1, 2 + $col).Value = $headers[$col]
}

$rows = @(
    @("050001","博时价值增长混合","21.87","74.33","2.75","0.6014",8),
    @("014611","富国核心科技12个月持有期混合A","5.97","81.00","5.44","0.3248",2),
    @("161219","国投瑞银新兴产业混合（LOF）","6.39","79.87","4.95","0.3163",5),
    @("050201","博时价值增长贰号混合","9.54","74.96","2.71","0.2585",8),
    @("016524","招商均衡成长混合A","3.47","80.67","3.86","0.1339",4),
    @("013630","嘉实均衡臻选一年持有期混合A","1.92","81.00","5.44","0.1044",2),
    @("014612","富国核心科技12个月持有期混合C","0.56","81.00","5.44","0.0305",2),
    @("519097","新华中小市值优选混合","0.66","70.51","4.36","0.0288",3),
    @("016525","招商均衡成长混合C","0.41","80.67","3.86","0.0158",4),
    @("013584","招商品质领航混合C","0.05","62.41","6.53","0.0033",1),
    @("013583","招商品质领航混合A","0.02","62.41","6.53","0.0013",1)
)

# Columns B, D, E, F, G must stay TEXT (fund codes / percentages with
# leading zeros & fixed decimals would otherwise be auto-coerced to
# numbers by Excel). Pre-format those columns as text before writing.
$newSheet.Range("B2:B12").NumberFormat = "@"
$newSheet.Range("D2:G12").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $row = $rows[$i]
    $newSheet.Cells.Item($r, 1).Value = $i          # A: zero-based index
    $newSheet.Cells.Item($r, 2).Value = $row[0]     # B: 基金代码 (text)
    $newSheet.Cells.Item($r, 3).Value = $row[1]     # C: 基金名称
    $newSheet.Cells.Item($r, 4).Value = $row[2]     # D: 基金规模 (text)
    $newSheet.Cells.Item($r, 5).Value = $row[3]     # E: 股票总仓位 (text)
    $newSheet.Cells.Item($r, 6).Value = $row[4]     # F: 仓位占比 (text)
    $newSheet.Cells.Item($r, 7).Value = $row[5]     # G: 持有市值(亿元) (text)
    $newSheet.Cells.Item($r, 8).Value = $row[6]     # H: 仓位排名 (number)
}

# --- 3) Move the fully-populated sheet to slot 2 (right after 总计) ------
$newSheet.Move($null, $totalSheet)

# --- 4) Update the "总计" summary sheet: insert a 2022-Q4 row on top ------
$totalSheet.Rows.Item(2).Insert() | Out-Null

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 11
$totalSheet.Cells.Item(2, 4).Value = 1.82

# Renumber the zero-based index column (A) for every row that shifted down
for ($r = 3; $r -le 10; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
